$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (replicationPolicy): add missing parentheses to CNReplication.setReplicationPolicy
$ws.Range("E10").Value = "CNReplication.setReplicationPolicy()"

# Row 9 (accessPolicy): CNAuthorization.systemMetadataChanged() -> CNCore.systemMetadataChanged()
$ws.Range("E9").Value = "manual (Tier 1), MNAuthorization.setAccessPolicy(), MNStorage.update ()(all must call CNCore.systemMetadataChanged())"

# Match the selection state recorded in the saved workbook
$ws.Range("E23").Select() | Out-Null
